$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.857.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.863.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5041"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3643"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07159"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8911"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07512"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.856.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.48%  "
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008508"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.902.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.018"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.094.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.404"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("E26").Value = "  -3.61%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.063"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.656"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09159"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05128"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7474"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.976"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.191"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.560"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01997"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5579"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.86%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.566"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.524"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.48%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4693"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.557"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.45%  "
